$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "maize silage" was renamed to "corn silage" in the substrate list (cell A2)
$ws.Range("A2").Value = "corn silage"

# The active selection in the sheet moved from C7 to D3
$ws.Range("D3").Select()

# Reflect the updated Excel window geometry (position/size) from the author's session
$win = $wb.Windows.Item(1)
$win.Left = 18800
$win.Top = 500
$win.Width = 10000
$win.Height = 15940
